$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new header column on sheet1 first, copying D1's formatting (bold,  ---
# --- centered, bordered) so E1 matches the other header cells exactly.          ---
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("E1").Value2 = "evaluator_partial_correctness"

# --- Create the two new sheets by copying sheet1 (carries the new header + styles) ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "o_20"

$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "o_20_jumbled"

# --- sheet1 (o_10): update llm_response / evaluator_response, add partial-correctness columns ---
$s1_llm = @'
Based on the given adjacency matrix, we can determine the shortest path from node A to node J by performing a breadth-first search (BFS) algorithm.
Starting from node A, we traverse the graph by visiting its neighboring nodes first. Then, we visit the neighbors of the visited nodes until we reach node J.
Here is the step-by-step process:
1. Start at node A.
2. Enqueue node A in a queue.
3. Set a boolean array to keep track of visited nodes and mark node A as visited.
4. Initialize a parent array to keep track of the parent nodes while traversing the graph.
5. While the queue is not empty:
   - Dequeue the first node from the queue.
   - If the dequeued node is equal to node J, stop the traversal.
   - Else, enqueue all the unvisited neighbors of the dequeued node and mark them as visited.
     - Also, set the parent of the newly visited nodes as the dequeued node.
6. Reconstruct the shortest path from node A to node J using the parent array.
Using this algorithm, we can find the shortest path from node A to node J, which is:
A -> B -> C -> D -> I -> J
'@
$ws1.Range("C2").Value2 = $s1_llm
$ws1.Range("D2").Value2 = "Wrong"
$s1_eval = @'
Output: 5/10
'@
$ws1.Range("E2").Value2 = $s1_eval

# --- sheet2 (o_20): new 20-node prompt/solution/response ---
$s2_prompt = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node X?
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T -> U -> V -> W -> X
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node T?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
'@
$ws2.Range("A2").Value2 = $s2_prompt

$s2_solution = @'
A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T
'@
$ws2.Range("B2").Value2 = $s2_solution

$s2_llm = @'
The shortest path from node A to node T is A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T.
'@
$ws2.Range("C2").Value2 = $s2_llm
$ws2.Range("D2").Value2 = "Correct"
$ws2.Range("E2").Value2 = "1. 20/20"

# --- sheet3 (o_20_jumbled): same prompt/solution, different response wording ---
$ws3.Range("A2").Value2 = $s2_prompt
$ws3.Range("B2").Value2 = $s2_solution

$s3_llm = @'
The shortest path from node A to node T is:
A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T
'@
$ws3.Range("C2").Value2 = $s3_llm
$ws3.Range("D2").Value2 = "Correct"
$ws3.Range("E2").Value2 = "Output: 20/20"

# --- Restore default (non-custom) row height on row 2 of every sheet; writing  ---
# --- the multi-line text above auto-expanded it, same as it would in real Excel. ---
$ws1.Rows.Item(2).AutoFit()
$ws2.Rows.Item(2).AutoFit()
$ws3.Rows.Item(2).AutoFit()

# --- Keep the originally active sheet (o_10) selected/active, as in the source file ---
$ws1.Activate()
